# Adds a new "InvalidLogin" worksheet (with sample invalid credentials)
# after the existing "ValidLogin" sheet, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# The existing (first) worksheet - "ValidLogin".
$wsValid = $wb.Worksheets.Item(1)

# Insert the new worksheet right after "ValidLogin" so it becomes sheet #2.
$wsInvalid = $wb.Worksheets.Add($null, $wsValid)
$wsInvalid.Name = "InvalidLogin"

# Same headers as the ValidLogin sheet, with invalid sample credentials.
$wsInvalid.Range("A1").Value = "Username"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"

# Match the saved selection/view state for each sheet.
$wsValid.Range("A1:B2").Select()
$wsInvalid.Range("B3").Select()

# Leave "InvalidLogin" as the active sheet/tab.
$wsInvalid.Select()
